$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "mejorada la pantalla de subida excel" -- add a new "RESENA" column to the
# upload template so reviewers can capture a short review/summary text.
#
# New header cell L1 = "RESENA". Copy K1's format first (PasteSpecial with
# formats only) so the new header gets the same bold/filled header style
# instead of picking up a brand-new style index, then set its text.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L1").Value = "RESENA"

# Reflect the new column in the saved view: scroll one column to the right
# (so column D becomes the left-most visible column) and move the active
# selection onto the new column's second row, just like the author did
# after adding the field.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L2").Select()
